# Article Type controller renamed to Post Type Controller and moved to admin
# module. Add new functionality for add privileges for roles. Update the
# "models" sheet: flip the Right/+ rows to Right/- and append four new
# resource sections (Country, Privilege, User, Post Type).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("models")

# Template ranges already present on the sheet, reused so the new cells pick
# up the same number formats / fills / fonts (Accent1 header, italic B-label,
# Good "+" / Bad "-" indicator) as their siblings instead of minting new
# style records.
$headerTemplate = $ws.Range("A8:K8")   # Accent1 section-header row
$itemBad        = $ws.Range("C9")      # "-" / Bad style sample
$itemGood       = $ws.Range("C6")      # "+" / Good style sample
$bGetItem       = $ws.Range("B9")      # "getItem" label sample
$bGetAll        = $ws.Range("B10")     # "getAll" label sample

function Add-Section {
    param($HeaderRow, $Title, $GetItemOk, $GetAllOk)

    $headerTemplate.Copy($ws.Range("A$HeaderRow`:K$HeaderRow"))
    $ws.Range("A$HeaderRow").Value2 = $Title

    $itemRow = $HeaderRow + 1
    $allRow  = $HeaderRow + 2

    $bGetItem.Copy($ws.Range("B$itemRow"))
    if ($GetItemOk) { $itemGood.Copy($ws.Range("C$itemRow")) } else { $itemBad.Copy($ws.Range("C$itemRow")) }

    $bGetAll.Copy($ws.Range("B$allRow"))
    if ($GetAllOk) { $itemGood.Copy($ws.Range("C$allRow")) } else { $itemBad.Copy($ws.Range("C$allRow")) }
}

# --- Existing "Right" section: both rows flip from "+" (Good) to "-" (Bad) ---
$itemBad.Copy($ws.Range("C3"))
$itemBad.Copy($ws.Range("C4"))

# --- New sections appended after the existing "Role" block (row 13) ---
Add-Section 14 "Country"   $false $false
Add-Section 17 "Privilege" $true  $true
Add-Section 20 "User"      $false $true
Add-Section 23 "Post Type" $true  $true

$ws.Range("C39").Select() | Out-Null
